$wb = $excel.ActiveWorkbook

# Sheet "Overview": Latest HO Xliff Generate Date (column G, row 2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-24 21:08:42"

# Sheet "zh-cn": Correspond Handoff Datetime (H2) and Correspond Handback DateTime (K2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-24 21:08:36"
$wsZhCn.Range("K2").Value = "2016-08-24 21:08:53"

# Sheet "de-de": Correspond Handback DateTime (K2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-24 21:09:03"
